$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add formula in B11: 10 times the value of C3
$ws.Range("B11").Formula = "=10*C3"

# Move the selection to B12, matching the post-edit cursor position
$ws.Range("B12").Select()
